# Append the latest Adafruit IO reading as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 21

# Column C ("Value") holds a numeric-looking reading ("25") but the feed
# logs every column as plain text, so force that cell to Text format first
# -- otherwise Excel's normal type inference would store it as a number.
$ws.Range("C" + $newRow).NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2024-09-25T18:06:40Z"
$ws.Range("B" + $newRow).Value = "temperature"
$ws.Range("C" + $newRow).Value = "25"
$ws.Range("D" + $newRow).Value = "N/A"
$ws.Range("E" + $newRow).Value = "N/A"
$ws.Range("F" + $newRow).Value = "N/A"
